# Refresh the cryptocurrency ranking table (Coin/Link/Price/Volume(1h))
# with the latest scrape, one row at a time. Cells whose new text still
# "looks like a number" (e.g. "7.07") are entered with a leading apostrophe
# so Excel keeps storing them as text (matching column D's existing cells),
# then the quote-prefix cell style picked up along the way is reset back to
# "Normal" so no cell formatting actually changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.420.54"
$ws.Range("E2").Value = "  +0.13%  "

# Row 3
$ws.Range("D3").Value = "3.662.82"
$ws.Range("E3").Value = "  -0.66%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").Value = "'645.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.20%  "

# Row 6
$ws.Range("D6").Value = "'159.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").Value = "'0.496"
$ws.Range("D8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = "  -0.91%  "

# Row 10
$ws.Range("D10").Value = "'7.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.25%  "

# Row 11
$ws.Range("D11").Value = "'0.438"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.23%  "

# Row 12
$ws.Range("D12").Value = "'0.0000230"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.56%  "

# Row 13
$ws.Range("D13").Value = "4.276.85"
$ws.Range("E13").Value = "  -0.70%  "

# Row 14
$ws.Range("D14").Value = "'32.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.82%  "

# Row 15
$ws.Range("D15").Value = "3.691.37"
$ws.Range("E15").Value = "  +0.15%  "

# Row 16
$ws.Range("D16").Value = "69.417.94"
$ws.Range("E16").Value = "  +0.15%  "

# Row 17
$ws.Range("E17").Value = "  +0.92%  "

# Row 18
$ws.Range("D18").Value = "'15.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.72%  "

# Row 19
$ws.Range("D19").Value = "'6.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.42%  "

# Row 20
$ws.Range("D20").Value = "'466.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.46%  "

# Row 21
$ws.Range("D21").Value = "'9.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.49%  "

# Row 22
$ws.Range("D22").Value = "'0.642"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.82%  "

# Row 23
$ws.Range("D23").Value = "'79.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.56%  "

# Row 24
$ws.Range("D24").Value = "3.807.03"
$ws.Range("E24").Value = "  -0.67%  "

# Row 25
$ws.Range("E25").Value = "  +0.07%  "

# Row 26
$ws.Range("D26").Value = "'0.0000124"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.85%  "

# Row 27
$ws.Range("D27").Value = "'10.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.21%  "

# Row 28
$ws.Range("D28").Value = "'8.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.01%  "

# Row 29
$ws.Range("E29").Value = "  -3.69%  "

# Row 30
$ws.Range("D30").Value = "'1.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.39%  "

# Row 31
$ws.Range("D31").Value = "'0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.17%  "

# Row 32
$ws.Range("D32").Value = "'1.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.76%  "

# Row 33
$ws.Range("D33").Value = "'26.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.42%  "

# Row 34
$ws.Range("D34").Value = "'6.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.54%  "

# Row 35
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.162"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.11%  "

# Row 36
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.651.10"
$ws.Range("E36").Value = "  -0.67%  "

# Row 37
$ws.Range("D37").Value = "'8.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.46%  "

# Row 39
$ws.Range("D39").Value = "'5.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.56%  "

# Row 40
$ws.Range("D40").Value = "'178.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.33%  "

# Row 41
$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.06%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'2.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.17%  "

# Row 43
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "'0.0888"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.26%  "

# Row 44
$ws.Range("D44").Value = "'0.926"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.83%  "

# Row 45
$ws.Range("D45").Value = "'46.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.73%  "

# Row 46
$ws.Range("D46").Value = "'2.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "

# Row 47
$ws.Range("E47").Value = "  -3.37%  "

# Row 48
$ws.Range("D48").Value = "'26.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.34%  "

# Row 49
$ws.Range("D49").Value = "'7.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.13%  "

# Row 50
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "'0.000264"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.30%  "

# Row 51
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").Value = "'1.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.99%  "
